# SUB1 DSMCRIT overlap.xlsx - break out "Unknown SUB1" into opioid subtypes
# (Heroin / Non-script methadone / Other opioids) as new column headers, and
# break out the DSMCRIT rows into Opioid Dep./Opioid Abuse/All Others/Unknown
# as a new leading row-category column. (Ran linear multivariate regression)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- structural edit: insert a new column A and a new row 1 -------------
$ws.Columns("A:A").Insert()
$ws.Rows("1:1").Insert()

# --- new row-category labels in column A (rows 3-6) ----------------------
$ws.Range("A3").Value = "Opioid Dep."
$ws.Range("A4").Value = "Opioid Abuse"
$ws.Range("A5").Value = "All Others"
$ws.Range("A6").Value = "Unknown"

# --- new second-level header row (row 1): opioid sub-breakdown -----------
$ws.Range("C1").Value = "Heroin"
$ws.Range("D1").Value = "Non-script methadone"
$ws.Range("E1").Value = "Other opioids"
$ws.Range("F1").Value = "All Others"
$ws.Range("G1").Value = "Unknown"

# wrap the two longer headers so the (taller) row still fits
$ws.Range("D1:E1").WrapText = $true
$ws.Rows("1:1").RowHeight = 34

# --- fill in a previously-blank data point --------------------------------
$ws.Range("E5").Value = 0

# --- formatting: bold the Subtotal row/column headers + totals -----------
$ws.Range("H2").Font.Bold = $true
$ws.Range("B7:H7").Font.Bold = $true
$ws.Range("C7:H7").Font.Bold = $true
$ws.Range("H3:H7").Font.Bold = $true

# --- borders: box around the data block C3:H7 -----------------------------
$top = $ws.Range("C3:H3").Borders.Item(8)
$top.LineStyle = 1
$top.Weight = 2

$bottom = $ws.Range("C7:H7").Borders.Item(9)
$bottom.LineStyle = 1
$bottom.Weight = 2

$left = $ws.Range("C3:C7").Borders.Item(7)
$left.LineStyle = 1
$left.Weight = 2

$right = $ws.Range("H3:H7").Borders.Item(10)
$right.LineStyle = 1
$right.Weight = 2

# --- window/view cosmetics -------------------------------------------------
$excel.ActiveWindow.Zoom = 150
$ws.Range("G12").Select()
